# Reorder "Recorded By" (column G) values so that a literal "System"
# token (case-sensitive, capital S) is moved to the front of the
# comma-separated list, keeping the remaining tokens in their original
# relative order. If no literal "System" token exists in the list, the
# tokens are simply reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $hasSystem = $true }
        }

        if ($hasSystem) {
            $rest = @()
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) { $rest += $p }
            }
            $newParts = @("System") + $rest
        } else {
            $newParts = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $newParts += $parts[$i]
            }
        }

        $newVal = $newParts -join ", "

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
